$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data to append below the existing table (row 94 -> row 95).
$row = 95

# Column A holds a date-formatted string ("2025/10/12"). Assigning it
# directly would make Excel auto-convert it into a date serial number,
# so we briefly force a text number format, set the value, then restore
# the cell to the default "Normal" style so no stray style index lingers.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2025/10/12"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").Value = "日"
$ws.Range("C$row").Value = 20
$ws.Range("D$row").Value = 38
